# Update the "K" column (column G) values on Sheet1.
# The commit regenerates the save_data sheet: column G (header "K", the
# strikeouts-per-outing figure, replacing the old "Strike#" stat) is
# recalculated/rewritten with new values, along with std/mean and s_vals
# elsewhere in the workbook. Here we just need to write the regenerated
# K values into column G for the affected data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    "2"  = 0
    "3"  = 2
    "4"  = 2
    "5"  = 4
    "6"  = 0
    "7"  = 1
    "8"  = 1
    "9"  = 1
    "10" = 1
    "11" = 1
    "12" = 0
    "13" = 2
    "14" = 1
    "15" = 1
    "16" = 1
    "17" = 0
    "18" = 1
    "19" = 0
    "20" = 0
    "21" = 1
    "22" = 2
    "23" = 1
    "24" = 1
    "25" = 2
    "26" = 1
    "27" = 1
    "28" = 2
    "29" = 1
    "30" = 0
    "31" = 2
    "32" = 0
    "33" = 2
    "34" = 0
    "36" = 2
    "37" = 2
    "38" = 1
    "39" = 2
    "40" = 1
    "41" = 0
    "42" = 0
    "43" = 3
    "44" = 0
    "45" = 0
    "46" = 1
    "47" = 0
    "48" = 2
    "49" = 1
    "50" = 2
    "51" = 3
    "52" = 2
    "53" = 2
    "54" = 1
    "55" = 2
    "56" = 1
    "57" = 1
    "58" = 1
    "60" = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
